# Campaign Global class as Reference + procedural teste generator
# Adds a new row (35 / SupplySquad / PassiveSkill) to the skills table,
# mirroring the formatting of the row above it, and leaves the new
# entry's name cell selected (matching the author's final cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New data row -----------------------------------------------------
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "SupplySquad"
$ws.Range("C37").Value = "PassiveSkill"
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0

# Carry over the same cell formatting used by the rest of the table
# (copy format only from the row directly above, like the other rows).
$ws.Range("A36:E36").Copy()
$ws.Range("A37:E37").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# PasteSpecial(xlPasteFormats) only touches formatting, but re-assert
# the values to be safe in case a host implementation clears them.
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = "SupplySquad"
$ws.Range("C37").Value = "PassiveSkill"
$ws.Range("D37").Value = 0
$ws.Range("E37").Value = 0

# --- View / selection state --------------------------------------------
# Scroll the sheet so the new row is visible, and leave the cursor on
# the new skill's name cell, matching where editing finished.
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B37").Select() | Out-Null
